$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A and B get an across-the-board offset applied (A += 0.2, B += 0.05)
# for rows 1-14, reflecting the "too great of a response" (overshoot) fix
# mentioned in the commit message. Row 14 (columns A-H) gets its own
# distinct updated values.
for ($r = 1; $r -le 13; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $cellA.Value() + 0.2

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = $cellB.Value() + 0.05
}

# Row 14 explicit new values
$ws.Cells.Item(14, 1).Value = 0.7
$ws.Cells.Item(14, 2).Value = 0.55
$ws.Cells.Item(14, 3).Value = 0.5
$ws.Cells.Item(14, 4).Value = 0.5
$ws.Cells.Item(14, 5).Value = 0.5
$ws.Cells.Item(14, 6).Value = 0.5
$ws.Cells.Item(14, 7).Value = 0.5
$ws.Cells.Item(14, 8).Value = 0.5

# Update the selection to match the saved view state
$ws.Range("L10").Select()
